$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @{
    10  = "32.179540,34.908235"
    34  = "32.707298,35.173307"
    49  = "32.808325,35.060120"
    50  = "31.748498,35.214655"
    56  = "32.163351,34.809456"
    88  = "32.064156,34.854185"
    117 = "31.942541,34.872538"
    128 = "31.753295,34.996429"
    130 = "31.749399,35.210830"
    132 = "31.857912,35.215438"
    139 = "31.225747,34.809580"
    140 = "31.068028,35.007787"
    141 = "31.238529,34.795441"
    144 = "31.863818,34.742477"
    145 = "31.928328,34.878378"
    146 = "32.175827,34.926297"
    147 = "32.045844,34.752383"
    153 = "31.665784,34.601137"
    157 = "31.419806,34.603236"
    163 = "32.093937,34.885592"
    164 = "31.244467,34.807280"
    166 = "31.757029,34.990864"
}

foreach ($row in $updates.Keys) {
    $ws.Range("B$row").Value = $updates[$row]
}
